# Plant Sim, presentation.pptx — apply author's edit via PowerPoint COM-interop
#
# Summary of the change (reconstructed from the OOXML diff):
#   1. Slide "Execution" (slide 3) gets extra bullet text split across runs,
#      two new bullet paragraphs inserted, and a superscript "th" in "4th".
#   2. Slide "Organization" moves from position 6 to position 4, gets a new
#      second paragraph of body text, and its picture shifts down to make
#      room for the new paragraph.
#   3. A brand-new slide "Problems & Soultions" is inserted right after
#      "Organization" (new position 5), using the same Title-and-body layout.
#   4. The rest of the deck (Web-page, Data Base, Resulting Prototype) keeps
#      its content and relative order, sliding down in slide position.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) "Execution" slide (position 3): rewrite the bulleted body text.
# ---------------------------------------------------------------------------
$execSlide = $p.Slides.Item(3)
$execBody  = $execSlide.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: "The app was made with web programing" -> longer sentence.
$execBody.Paragraphs(1).Text = "-The app was made with web programing, a web server, and a SQL database, and the back-end code was done in Javascript"

# Paragraph 2: "Files were exchanged with the group through github. " -> trimmed trailing period wording.
$execBody.Paragraphs(2).Text = "- Files were exchanged with the group through github"

# Insert two brand-new bullet paragraphs right after paragraph 2 (before the
# "https://github.com/..." paragraph).
$execBody.Paragraphs(2).InsertAfter("`r- Communitiation was done with im, weekly meetings, and sprints `r- A prototype was completed by the end of the 4th sprint. ") | Out-Null

# Make the "th" in "4th" a superscript, matching the author's formatting.
$newPara = $execBody.Paragraphs(4)
$text4 = $newPara.Text
$thPos = $text4.IndexOf("4th") + 2   # 1-based position of "t" in "th"
$thRange = $newPara.Characters($thPos, 2)
$thRange.Font.Superscript = $true

# ---------------------------------------------------------------------------
# 2) Move "Organization" slide from position 6 to position 4, add a second
#    body paragraph, and push the picture down to make room for it.
# ---------------------------------------------------------------------------
$orgSlide = $p.Slides.Item(6)
$orgSlide.MoveTo(4)

$orgSlide = $p.Slides.Item(4)
$orgBody  = $orgSlide.Shapes.Item(2).TextFrame.TextRange
$orgBody.Text = "The group worked in sprints, and utilized agile methods we learned. Sample product backlog planning. `rHere" + [char]8217 + "s an example of our sprint backlog"

$orgPic = $orgSlide.Shapes.Item(3)
$orgPic.Top = 2615295 / 12700

# ---------------------------------------------------------------------------
# 3) Insert the brand-new "Problems & Soultions" slide at position 5.
# ---------------------------------------------------------------------------
$newLayout = $orgSlide.CustomLayout
$probSlide = $p.Slides.AddSlide(5, $newLayout)
$probSlide.Shapes.Item(1).Name = "Title 1"
$probSlide.Shapes.Item(2).Name = "Text Placeholder 2"

$probSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Problems & Soultions"

$probBody = $probSlide.Shapes.Item(2).TextFrame.TextRange
$probBody.Text = "One of the difficult things about this project that came up in the sprints consitiaintly was the integration between backend code, database and gui."

$probParas = @(
    "Since everyone knew Java we started with that, but had difficulty finding a good way to integrate it effiecntly in the web. During our incremental development through sprints, we switch to java script for something more fesiable. ",
    "",
    "",
    "Scrum is one of the widely used Agile Methodologies. In Scrum we do incremental and iterative development and these iterations are termed as Sprints. These Sprints are usually time boxed to 2-4 weeks.",
    "Before each Sprint starts, team decides which functionality or user stories (a software system feature specified by the customer) will be incorporated and developed during this Sprint.",
    ""
)
foreach ($para in $probParas) {
    $probBody.InsertAfter("`r" + $para) | Out-Null
}
